# Fruta / hortaliza, semanal
#
# Weekly refresh of the "Mango" sheet: a brand-new week's record is
# inserted at the top of the data block (row 121), pushing every
# existing data row down by one (old row 121 -> 122, ... old row 134
# -> 135).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data block (rows 121:134) down by one row, carrying
# formatting (e.g. the date-formatted style on column D) along with it.
$ws.Rows("121:121").Insert()

# Populate the newly opened row 121 with this week's record.
$ws.Range("A121").Value = 5
$ws.Range("B121").Value = 'Macroferia Regional de Talca'
$ws.Range("C121").Value = 'Maule'
$ws.Range("D121").Value = 44748
$ws.Range("E121").Value = 7
$ws.Range("F121").Value = 'Fruta'
$ws.Range("G121").Value = 100108
$ws.Range("H121").Value = 'Tropicales y subtropicales'
$ws.Range("I121").Value = 100108002
$ws.Range("J121").Value = 'Mango'
$ws.Range("K121").Value = 'Sin especificar'
$ws.Range("L121").Value = 'Primera'
$ws.Range("M121").Value = 240
$ws.Range("N121").Value = 8000
$ws.Range("O121").Value = 8000
$ws.Range("P121").Value = 8000
$ws.Range("Q121").Value = '$/bandeja 4 kilos'
$ws.Range("R121").Value = 'Brasil'
$ws.Range("S121").Value = 2000
$ws.Range("T121").Value = 4
